# edit.ps1 -- apply the "Add files via upload" revision to Member-Contributions.docx
#
# Summary of changes:
#  1. "Each member's ... So if your team ..." paragraph: split the long run into
#     three runs, changing "So if" -> "So, if" (the comma is typed and the
#     resulting "So," becomes its own run).
#  2. Table header cell "% of Contribution " is split into three runs around the
#     word "of", which gets wrapped in gramStart/gramEnd proofErr markers.
#  3. Table cell "Chisato Sakata" is split into two runs ("Chisato" / " Sakata"),
#     with spellStart/spellEnd proofErr markers around "Chisato".
#  4. Table cell "A.T" gets wrapped in gramStart/gramEnd proofErr markers.
#  5. A new bullet item "Helped with additional research on topic" is added
#     after "Worked on the PowerPoint." in Middy Esmail's comments cell.
#
# All edits are applied by replacing the OOXML of the whole affected
# paragraph(s) with an exact, hand-built paragraph (via Range.InsertXML),
# which lets us reproduce the run-splits and <w:proofErr/> markers exactly
# the way Word itself would leave them after typing + a spelling/grammar
# pass, rather than relying on incidental run-merging behaviour of plain
# text insertion.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "Each member's final grade ... So if your team ..." paragraph
# ---------------------------------------------------------------------
$r1 = $d.Content
$r1.Find.MatchCase = $true
$r1.Find.Execute("Each member")
if (-not $r1.Find.Found) { throw "paragraph 1 anchor text not found" }
$para1 = $r1.Paragraphs(1).Range

$xml1 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="75056A4C" w14:textId="44457581" w:rsidR="0055546D" w:rsidRDefault="0055546D" w:rsidP="0055546D"><w:pPr><w:rPr><w:rFonts w:ascii="Cambria" w:hAnsi="Cambria" w:cstheme="minorHAnsi"/></w:rPr></w:pPr><w:r w:rsidRPr="0055546D"><w:rPr><w:rFonts w:ascii="Cambria" w:hAnsi="Cambria" w:cstheme="minorHAnsi"/></w:rPr><w:t xml:space="preserve">Each member’s final grade on the assignment will be their % on this form of your team’s overall grade.  </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Cambria" w:hAnsi="Cambria" w:cstheme="minorHAnsi"/></w:rPr><w:t>So,</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Cambria" w:hAnsi="Cambria" w:cstheme="minorHAnsi"/></w:rPr><w:t xml:space="preserve"> if your team gets 80% for your tutorial, a member making a half contribution (50% on this form) will get a 40% on the assignment.</w:t></w:r><w:r w:rsidR="00D1149B"><w:rPr><w:rFonts w:ascii="Cambria" w:hAnsi="Cambria" w:cstheme="minorHAnsi"/></w:rPr><w:t xml:space="preserve">  A member making a 100% contribution will get an 80%.</w:t></w:r></w:p>
'@
$para1.InsertXML($xml1)

# ---------------------------------------------------------------------
# 2) Table header cell: "% of Contribution " -> "% " / "of" / " Contribution "
# ---------------------------------------------------------------------
$r2 = $d.Content
$r2.Find.Execute("% of Contribution")
if (-not $r2.Find.Found) { throw "paragraph 2 anchor text not found" }
$para2 = $r2.Paragraphs(1).Range

$xml2 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="61312787" w14:textId="508D9415" w:rsidR="00A7085F" w:rsidRPr="0055546D" w:rsidRDefault="00A7085F" w:rsidP="00BB7FB6"><w:pPr><w:rPr><w:rFonts w:ascii="Cambria" w:hAnsi="Cambria" w:cstheme="majorHAnsi"/></w:rPr></w:pPr><w:r w:rsidRPr="0055546D"><w:rPr><w:rFonts w:ascii="Cambria" w:hAnsi="Cambria" w:cstheme="majorHAnsi"/></w:rPr><w:t xml:space="preserve">% </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Cambria" w:hAnsi="Cambria" w:cstheme="majorHAnsi"/></w:rPr><w:t>of</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="Cambria" w:hAnsi="Cambria" w:cstheme="majorHAnsi"/></w:rPr><w:t xml:space="preserve"> Contribution </w:t></w:r></w:p>
'@
$para2.InsertXML($xml2)

# ---------------------------------------------------------------------
# 3) Table cell: "Chisato Sakata" -> "Chisato" (spellStart/spellEnd) + " Sakata"
# ---------------------------------------------------------------------
$r3 = $d.Content
$r3.Find.Execute("Chisato Sakata")
if (-not $r3.Find.Found) { throw "paragraph 3 anchor text not found" }
$para3 = $r3.Paragraphs(1).Range

$xml3 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="26202B88" w14:textId="17E08979" w:rsidR="00A7085F" w:rsidRPr="0055546D" w:rsidRDefault="00AA4B4A"><w:pPr><w:rPr><w:rFonts w:ascii="Cambria" w:hAnsi="Cambria" w:cstheme="majorHAnsi"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00AA4B4A"><w:rPr><w:rFonts w:ascii="Cambria" w:hAnsi="Cambria" w:cstheme="majorHAnsi"/></w:rPr><w:t>Chisato</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Cambria" w:hAnsi="Cambria" w:cstheme="majorHAnsi"/></w:rPr><w:t xml:space="preserve"> Sakata</w:t></w:r></w:p>
'@
$para3.InsertXML($xml3)

# ---------------------------------------------------------------------
# 4) Table cell: "A.T" gets wrapped with gramStart/gramEnd proofErr markers
# ---------------------------------------------------------------------
$r4 = $d.Content
$r4.Find.Execute("A.T")
if (-not $r4.Find.Found) { throw "paragraph 4 anchor text not found" }
$para4 = $r4.Paragraphs(1).Range

$xml4 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="3BB48122" w14:textId="79D3DAAC" w:rsidR="00A7085F" w:rsidRPr="0055546D" w:rsidRDefault="00FB3902"><w:pPr><w:rPr><w:rFonts w:ascii="Cambria" w:hAnsi="Cambria" w:cstheme="majorHAnsi"/></w:rPr></w:pPr><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Cambria" w:hAnsi="Cambria" w:cstheme="majorHAnsi"/></w:rPr><w:t>A.T</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p>
'@
$para4.InsertXML($xml4)

# ---------------------------------------------------------------------
# 5) Add new bullet "Helped with additional research on topic" after
#    "Worked on the PowerPoint." in the comments cell.
# ---------------------------------------------------------------------
$r5 = $d.Content
$r5.Find.Execute("Worked on the PowerPoint")
if (-not $r5.Find.Found) { throw "paragraph 5 anchor text not found" }
$para5 = $r5.Paragraphs(1).Range

$xml5 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="4EEB23F5" w14:textId="027B6058" w:rsidR="006B15D3" w:rsidRPr="006B15D3" w:rsidRDefault="006B15D3" w:rsidP="006B15D3"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:rFonts w:ascii="Cambria" w:hAnsi="Cambria" w:cstheme="majorHAnsi"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Cambria" w:hAnsi="Cambria" w:cstheme="majorHAnsi"/></w:rPr><w:t xml:space="preserve">Worked on the PowerPoint. </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:rFonts w:ascii="Cambria" w:hAnsi="Cambria" w:cstheme="majorHAnsi"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Cambria" w:hAnsi="Cambria" w:cstheme="majorHAnsi"/></w:rPr><w:t>Helped with additional research on topic</w:t></w:r></w:p>
'@
$para5.InsertXML($xml5)

Write-Output "All edits applied."
